$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 280 (pushing existing rows 280+ down by 2),
# so that "19_02_05_others" and "19_02_17_electricity" are added right
# after "19_02_04_biomass" (row 279) and before the "fuel_cell_ev" block.
$ws.Rows.Item(280).Resize(2).Insert()

$ws.Range("B280").Value = "19_02_05_others"
$ws.Range("C280").Value = "19_02_05_others"
$ws.Range("D280").Value = "(new)"

$ws.Range("B281").Value = "19_02_17_electricity"
$ws.Range("C281").Value = "19_02_17_electricity"
$ws.Range("D281").Value = "(new)"

# Reflect the author's on-screen selection after the edit.
$ws.Range("D279:D281").Select()
